$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Move the report title from D2 to B2 (Cut keeps formatting + value, clears source cell)
$ws.Range("D2").Cut($ws.Range("B2")) | Out-Null
$ws.Range("D2").Clear() | Out-Null

# 2. Fill in the previously-empty cells C4, B5, C5, C7 using the same format as
#    their already-styled neighbours (B4 / B7) so they pick up style index 1.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null
$ws.Range("B5:C5").PasteSpecial(-4122) | Out-Null

$ws.Range("B7").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null

# 3. Row 6: un-merge B6:E6, drop the old border/font style, reuse B7's plain style
#    for B6, and clear out C6:E6 entirely (no longer part of the layout).
$ws.Range("B6:E6").UnMerge() | Out-Null
$ws.Range("B7").Copy() | Out-Null
$ws.Range("B6").PasteSpecial(-4122) | Out-Null
$ws.Range("C6:E6").Clear() | Out-Null

# 4. Row 8 (Doanh nghiep / San pham / So luong truy cap header row): centre the text
#    -> turns the old border-only style into border+center-aligned style.
$ws.Range("B8:D8").HorizontalAlignment = -4108

# 5. Row 9 is no longer part of the report - remove it completely.
$ws.Range("B9:D9").EntireRow.Delete() | Out-Null

# 6. Column sizing to match the new, wider layout (values chosen so Excel's
#    pixel-snapped ColumnWidth lands on 49 / ~9.66 / ~19.78 characters).
$ws.Columns(2).ColumnWidth = 48.1667
$ws.Columns(3).ColumnWidth = 8.8333
$ws.Columns(4).ColumnWidth = 19.0

# 7. Selection cursor ends on D6 (matches the authored selection state).
$ws.Range("D6").Select() | Out-Null
